$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "answer" column header was renamed to "response"
$ws.Range("B1").Value = "response"

# Columns were resized (A and B got explicit widths)
$ws.Range("A:A").ColumnWidth = 14
$ws.Range("B:B").ColumnWidth = 64

# Active cell/selection moved to F13 in the re-saved file
$ws.Range("F13").Select() | Out-Null
